# Tasks for Code Reviews.xlsx - update for code review 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task descriptions (row 34 "deployment" task, row 35 "PopRep" task).
# Set C35 before C34 so the shared-string table entries land in the same
# order as the target workbook (index 29 = PopRep text, index 30 = deployment text).
$ws.Range("C35").Value = '"PopRep" aka one more report set up'
$ws.Range("C34").Value = "Supposed to do deployoment"

# Update percentages for code review 4 rows
$ws.Range("D33").Value = 0.85
$ws.Range("D34").Value = 0.05
$ws.Range("D35").Value = 0.1

# Update view state: zoom to 81% and move selection to G32
$excel.ActiveWindow.Zoom = 81
$ws.Range("G32").Select()
